$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Delete()
